$wb = $excel.ActiveWorkbook

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 2690949.2
$ws.Range("I138").Value2 = 969.1053000000001
$ws.Range("J138").Value2 = 3879545.2
$ws.Range("K138").Value2 = 2907.3159
$ws.Range("L138").Value2 = 11638635.6
$ws.Range("M138").Value2 = 2232.6841
$ws.Range("N138").Value2 = -11648915.6

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3708.12
$ws.Range("I32").Value2 = 2666.7952
$ws.Range("J32").Value2 = 8792.235000000001
$ws.Range("K32").Value2 = 2666.7952
$ws.Range("L32").Value2 = 8792.235000000001
$ws.Range("M32").Value2 = -2379.7952
$ws.Range("N32").Value2 = -9366.235000000001

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 83501210
$ws.Range("I61").Value2 = 111223050
$ws.Range("J61").Value2 = 335671.34
$ws.Range("K61").Value2 = 111223050
$ws.Range("L61").Value2 = 335671.34
$ws.Range("M61").Value2 = -111222838
$ws.Range("N61").Value2 = -336095.34

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 25203268
$ws.Range("I74").Value2 = 50201796
$ws.Range("J74").Value2 = 204740
$ws.Range("K74").Value2 = 50201796
$ws.Range("L74").Value2 = 204740
$ws.Range("M74").Value2 = -50200922
$ws.Range("N74").Value2 = -206488

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 25203268
$ws.Range("I77").Value2 = 50201796
$ws.Range("J77").Value2 = 204740
$ws.Range("K77").Value2 = 251008980
$ws.Range("L77").Value2 = 1023700
$ws.Range("M77").Value2 = -251004612
$ws.Range("N77").Value2 = -1032436

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 83501210
$ws.Range("I136").Value2 = 111223050
$ws.Range("J136").Value2 = 335671.34
$ws.Range("K136").Value2 = 333669150
$ws.Range("L136").Value2 = 1007014.02
$ws.Range("M136").Value2 = -333666600
$ws.Range("N136").Value2 = -1012114.02

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 3510.0256
$ws.Range("I134").Value2 = 2674.7856
$ws.Range("J134").Value2 = 5636.091
$ws.Range("K134").Value2 = 8024.3568
$ws.Range("L134").Value2 = 16908.273
$ws.Range("M134").Value2 = -5489.3568
$ws.Range("N134").Value2 = -21978.273

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 783461.5600000001
$ws.Range("I4").Value2 = 5000
$ws.Range("J4").Value2 = 1691666.6
$ws.Range("K4").Value2 = 5000
$ws.Range("L4").Value2 = 1691666.6
$ws.Range("M4").Value2 = -4888
$ws.Range("N4").Value2 = -1691890.6

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 266.25
$ws.Range("I7").Value2 = 75
$ws.Range("J7").Value2 = 330
$ws.Range("K7").Value2 = 75
$ws.Range("L7").Value2 = 330
$ws.Range("M7").Value2 = 38
$ws.Range("N7").Value2 = -556

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2560.8096
$ws.Range("I31").Value2 = 1388.4193
$ws.Range("J31").Value2 = 5864.8184
$ws.Range("K31").Value2 = 1388.4193
$ws.Range("L31").Value2 = 5864.8184
$ws.Range("M31").Value2 = -1093.4193
$ws.Range("N31").Value2 = -6454.8184

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 2560.8096
$ws.Range("I34").Value2 = 1388.4193
$ws.Range("J34").Value2 = 5864.8184
$ws.Range("K34").Value2 = 1388.4193
$ws.Range("L34").Value2 = 5864.8184
$ws.Range("M34").Value2 = -1186.4193
$ws.Range("N34").Value2 = -6268.8184

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 26317700
$ws.Range("I58").Value2 = 32259788
$ws.Range("K58").Value2 = 32259788
$ws.Range("M58").Value2 = -32259585

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 26317700
$ws.Range("I136").Value2 = 32259788
$ws.Range("K136").Value2 = 96779364
$ws.Range("M136").Value2 = -96776814

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 1088.566
$ws.Range("I131").Value2 = 461.14285
$ws.Range("J131").Value2 = 1184.0435
$ws.Range("K131").Value2 = 1383.42855
$ws.Range("L131").Value2 = 3552.1305
$ws.Range("M131").Value2 = 3656.57145
$ws.Range("N131").Value2 = -13632.1305

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value2 = 848.2353000000001
$ws.Range("I132").Value2 = 672.8570999999999
$ws.Range("J132").Value2 = 1666.6666
$ws.Range("K132").Value2 = 6055.7139
$ws.Range("L132").Value2 = 14999.9994
$ws.Range("M132").Value2 = -3525.7139
$ws.Range("N132").Value2 = -20059.9994

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value2 = 25477.24
$ws.Range("I137").Value2 = 1015.8
$ws.Range("J137").Value2 = 41784.867
$ws.Range("K137").Value2 = 3047.4
$ws.Range("L137").Value2 = 125354.601
$ws.Range("M137").Value2 = 2052.6
$ws.Range("N137").Value2 = -135554.601

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value2 = 2869.6155
$ws.Range("I138").Value2 = 1726.4286
$ws.Range("J138").Value2 = 4203.3335
$ws.Range("K138").Value2 = 5179.2858
$ws.Range("L138").Value2 = 12610.0005
$ws.Range("M138").Value2 = -39.28579999999965
$ws.Range("N138").Value2 = -22890.0005

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 4728.5713
$ws.Range("I5").Value2 = 2620
$ws.Range("J5").Value2 = 10000
$ws.Range("K5").Value2 = 2620
$ws.Range("L5").Value2 = 10000
$ws.Range("M5").Value2 = -2508
$ws.Range("N5").Value2 = -10224

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 4022.9333
$ws.Range("I80").Value2 = 3780
$ws.Range("J80").Value2 = 4040.2856
$ws.Range("K80").Value2 = 3780
$ws.Range("L80").Value2 = 4040.2856
$ws.Range("M80").Value2 = -2782
$ws.Range("N80").Value2 = -6036.2856

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value2 = 4022.9333
$ws.Range("I83").Value2 = 3780
$ws.Range("J83").Value2 = 4040.2856
$ws.Range("K83").Value2 = 18900
$ws.Range("L83").Value2 = 20201.428
$ws.Range("M83").Value2 = -13908
$ws.Range("N83").Value2 = -30185.428

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 1500
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 1500
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 1500
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value2 = -1724

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 240.57143
$ws.Range("I55").Value2 = 136
$ws.Range("J55").Value2 = 502
$ws.Range("K55").Value2 = 136
$ws.Range("L55").Value2 = 502
$ws.Range("M55").Value2 = 37
$ws.Range("N55").Value2 = -848

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 42165.527
$ws.Range("I132").Value2 = 17614.908
$ws.Range("K132").Value2 = 52844.724
$ws.Range("M132").Value2 = -50314.724

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value2 = 350766.66
$ws.Range("I2").Value2 = 619780
$ws.Range("J2").Value2 = 14500
$ws.Range("K2").Value2 = 619780
$ws.Range("L2").Value2 = 14500
$ws.Range("M2").Value2 = -619668
$ws.Range("N2").Value2 = -14724

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 4213.857
$ws.Range("I62").Value2 = 3664.6667
$ws.Range("K62").Value2 = 3664.6667
$ws.Range("M62").Value2 = -3040.6667

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value2 = 4213.857
$ws.Range("I65").Value2 = 3664.6667
$ws.Range("K65").Value2 = 18323.3335
$ws.Range("M65").Value2 = -15203.3335

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 97111.28999999999
$ws.Range("I132").Value2 = 78256.84
$ws.Range("K132").Value2 = 234770.52
$ws.Range("M132").Value2 = -232240.52
